$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.010.97'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '2.407.04'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '554.19'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '135.50'
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.583'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").Value = '5.62'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").Value = '24.59'
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").Value = '2.838.13'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '59.892.92'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").Value = '2.355.86'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").Value = '11.17'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").Value = '4.50'
$ws.Range("E19").Value = '  +3.28%  '
$ws.Range("D20").Value = '327.31'
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '64.61'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").Value = '0.180'
$ws.Range("E24").Value = '  +4.33%  '
$ws.Range("D25").Value = '8.60'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +4.49%  '
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").Value = '0.0₃0765'
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '169.46'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").Value = '6.13'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  +8.91%  '
$ws.Range("D33").Value = '0.400'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").Value = '18.42'
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = '1.32'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '4.18'
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").Value = '323.38'
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").Value = '147.01'
$ws.Range("E41").Value = '  +6.37%  '
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("D43").Value = '0.0964'
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").Value = '19.84'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").Value = '0.0514'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '0.575'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("E51").Value = '  -1.16%  '
